# "Update registro de empresas con Excel"
# Mark every existing company row in the "Registro de empresas" sheet as
# "do not modify" by filling the new flag column (I) with "X" for each
# data row (I2:I26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registro de empresas")

$ws.Range("I2:I26").Value = "X"
